# Generate Report for Handback
# Updates the timestamp values recorded on the "Overview", "zh-cn" and
# "de-de" worksheets to reflect a newly generated handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 - Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2017-02-21 10:39:10"

# zh-cn!H2 - Correspond Handoff Datetime
$wsZhCn.Range("H2").Value = "2017-02-21 10:38:53"

# zh-cn!L2 - Correspond Handback DateTime
$wsZhCn.Range("L2").Value = "2017-02-21 10:39:49"

# de-de!H2 - Correspond Handoff Datetime
$wsDeDe.Range("H2").Value = "2017-02-21 10:39:10"

# de-de!L2 - Correspond Handback DateTime
$wsDeDe.Range("L2").Value = "2017-02-21 10:40:13"
